# Update TPM-derived values on Sheet1 to reflect the new expression table
# (commit: "update scripts wuth new tpm").
# Only the Receptor/Edge expression + specificity columns (M:T) on rows 2-6
# change; everything else (ids, counts, detection rates) is untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 9.409854000000001
$ws.Range("N2").Value = 28.229562
$ws.Range("O2").Value = 0.151097812489362
$ws.Range("P2").Value = 0.1551735891706745
$ws.Range("Q2").Value = 15.122237608656
$ws.Range("R2").Value = 136.100138477904
$ws.Range("S2").Value = 0.151097812489362
$ws.Range("T2").Value = 0.1551735891706745

# Row 3 (only the derived-specificity columns changed)
$ws.Range("O3").Value = 0.2814274088726355
$ws.Range("P3").Value = 0.2890187515378117
$ws.Range("S3").Value = 0.2814274088726355
$ws.Range("T3").Value = 0.2890187515378117

# Row 4
$ws.Range("M4").Value = 16.72420333333333
$ws.Range("N4").Value = 50.17260999999999
$ws.Range("O4").Value = 0.2685472632512643
$ws.Range("P4").Value = 0.275791171388365
$ws.Range("Q4").Value = 26.87686510568
$ws.Range("R4").Value = 241.8917859511199
$ws.Range("S4").Value = 0.2685472632512643
$ws.Range("T4").Value = 0.275791171388365

# Row 5
$ws.Range("M5").Value = 4.907254
$ws.Range("N5").Value = 9.814508
$ws.Range("O5").Value = 0.07879775230621766
$ws.Range("P5").Value = 0.05394885093521105
$ws.Range("Q5").Value = 7.886271242256
$ws.Range("R5").Value = 47.317627453536
$ws.Range("S5").Value = 0.07879775230621766
$ws.Range("T5").Value = 0.05394885093521105

# Row 6
$ws.Range("M6").Value = 13.70892733333334
$ws.Range("N6").Value = 41.12678200000001
$ws.Range("O6").Value = 0.2201297630805206
$ws.Range("P6").Value = 0.2260676369679378
$ws.Range("Q6").Value = 22.03112359601601
$ws.Range("R6").Value = 198.280112364144
$ws.Range("S6").Value = 0.2201297630805206
$ws.Range("T6").Value = 0.2260676369679378
